$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the original (pre-edit) values for the columns that change,
# for every data row (2-29), reading via Value2 (Value getter is unreliable
# for strings in this host).
$cols = @("D","K","L","M","N","O","P","Q","R","S","T")
$orig = @{}
for ($r = 2; $r -le 29; $r++) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $orig[$r] = $rowData
}

# Row permutation: target row -> source row (weekly re-pull reshuffled the
# daily price records across the existing date range).
$map = @{ 2=5; 3=6; 4=25; 5=11; 6=16; 7=23; 8=18; 9=3; 10=4; 11=17; 12=12; 13=14; 14=27; 15=28; 16=9; 17=10; 18=2; 19=20; 20=21; 21=22; 22=29; 23=7; 24=26; 25=19; 26=24; 27=8; 28=13; 29=15 }

foreach ($target in ($map.Keys | Sort-Object)) {
    $source = $map[$target]
    $data = $orig[$source]
    foreach ($c in $cols) {
        $ws.Range("$c$target").Value = $data[$c]
    }
}
